# Rigori_Ianesi.xlsx -- finish analysis of penalty takers for Pontedera
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet ---
$ws.Name = "Foglio1"

# --- 2. Correct the minute values for the two existing rows ---
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 27

# --- 3. Add the (empty) 8th column header, matching the other header cells ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").ClearContents()

# --- 4. Extend rows 2-3 formatting (font + fill) down through row 16 first,
#         so every new cell already carries the data font (fontId 1) before
#         any alignment / number-format tweak is layered on top of it. ---
$ws.Range("A2:H3").Copy()
$ws.Range("A4:H16").PasteSpecial(-4122)
$ws.Range("A4:H16").ClearContents()

# --- 5. Refresh the data-row font so it no longer inherits the theme scheme ---
$ws.Range("A2:H16").Font.Name = "Arial"

# --- 6. Right-align the numeric columns for every data row ---
$ws.Range("B2:B16").HorizontalAlignment = -4152
$ws.Range("C2:C16").HorizontalAlignment = -4152
$ws.Range("D2:D16").HorizontalAlignment = -4152
$ws.Range("E2:E16").HorizontalAlignment = -4152
$ws.Range("H2:H16").HorizontalAlignment = -4152

# --- 7. Mark score column as text so values like "0-0" are not reinterpreted ---
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F16").NumberFormat = "@"

# --- 8. Re-color + rename the hyperlink font used by the Link column ---
$ws.Range("G2:G16").Font.Color = 13391121
$ws.Range("G2:G16").Font.Name = "Arial"

# --- 9. Column A: drop the custom width back to the sheet default ---
$ws.Columns.Item(1).ColumnWidth = 12.63
